$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Driver class section: update points earned and grading comments
$ws.Range("E29").Value = 12
$ws.Range("F29").Value = "(-2) for not declaring and initailizing customer object, (-1) for passing incorrect arguments to product object, (-1) for not adding them to the inventory."
$ws.Range("F30").Value = "(-4) For no output displayed due to compilation errors"

# Generic section: update compilation errors grading comment
$ws.Range("F37").Value = "(-5) For compilation errors in Driver as well as other classes"

# Move the active selection to F37, matching the final state of the sheet
$ws.Range("F37").Select()

$wb.Save()
